# "updated commodities and sets"
#
# The "products" sheet previously split the "Neo bulk" product class into
# "Neo bulk (fast)" / "Neo bulk (slow)" product groups while "Break bulk"
# was a single (non-split) product group/class. This edit flips that: the
# "Break bulk" product class is now split into "Break bulk (fast)" /
# "Break bulk (slow)" product groups, while "Neo bulk" becomes a single
# (non-split) product group/class again.

$wb = $excel.ActiveWorkbook

$products = $wb.Worksheets.Item("products")

# Row 6: "Break bulk" -> "Break bulk (fast)" (product class column C6 stays "Break bulk")
$products.Range("B6").Value = "Break bulk (fast)"

# Row 7: "Neo bulk (fast)" / "Neo bulk" -> "Break bulk (slow)" / "Break bulk"
$products.Range("B7").Value = "Break bulk (slow)"
$products.Range("C7").Value = "Break bulk"

# Row 8: "Neo bulk (slow)" -> "Neo bulk" (product class column C8 stays "Neo bulk")
$products.Range("B8").Value = "Neo bulk"

# Reflect the author's final on-screen state: the "fuels" tab's zoom was
# reset back to 100% and the "products" sheet became the active tab with
# cell C8 selected.
$fuels = $wb.Worksheets.Item("fuels")
$fuels.Activate()
$excel.ActiveWindow.Zoom = 100

$products.Activate()
$products.Range("C8").Select()
